$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents (A1:K131) so the shared-strings table
# is rebuilt from scratch in the exact order we populate cells below.
$ws.Range("A1:K131").ClearContents()

# Re-populate cells in shared-string index order (0..66) so the resulting
# sharedStrings.xml table matches the target ordering.
$ws.Range("A1").Value = "חותמת זמן"
$ws.Range("C1").Value = "איך קוראים לך?"
$ws.Range("C2").Value = "זיגי"
$ws.Range("C3").Value = "טורצקי"
$ws.Range("C4").Value = "גלי"
$ws.Range("C5").Value = "הרצברג"
$ws.Range("C7").Value = "עמרי נצן"
$ws.Range("C8").Value = "שחר סלע"
$ws.Range("C9").Value = "פייביש"
$ws.Range("C10").Value = "דוידזון"
$ws.Range("C11").Value = "יניב"
$ws.Range("C12").Value = "מודן"
$ws.Range("C13").Value = "בן שושן"
$ws.Range("C14").Value = "זריהן"
$ws.Range("C15").Value = "מנחה"
$ws.Range("C16").Value = "שקד"
$ws.Range("C17").Value = "סוירי"
$ws.Range("C18").Value = "איילה"
$ws.Range("C19").Value = "פינסלר"
$ws.Range("C20").Value = "יובל בר"
$ws.Range("C21").Value = "ירין"
$ws.Range("C22").Value = "מיכאל"
$ws.Range("C23").Value = "עדי טל"
$ws.Range("C24").Value = "אורי ככ"
$ws.Range("C25").Value = "שושני"
$ws.Range("C26").Value = "אפרימה"
$ws.Range("C27").Value = "אסף בן חיים"
$ws.Range("C28").Value = "דור פרידמן"
$ws.Range("C29").Value = "יזהר"
$ws.Range("C30").Value = "ליעם מגיד"
$ws.Range("C6").Value = "ליאם מור"
$ws.Range("C31").Value = "נימי"
$ws.Range("C32").Value = "גליקמן"
$ws.Range("C33").Value = "עמית בר"
$ws.Range("C34").Value = "עמית לוי"
$ws.Range("C35").Value = "אליאב"
$ws.Range("C36").Value = "אריאל בן אליעזר"
$ws.Range("C37").Value = "דוד פיי"
$ws.Range("C38").Value = "יהל שפי"
$ws.Range("C39").Value = "יהלי מערבי ברנר"
$ws.Range("C40").Value = "אברג'ל"
$ws.Range("C41").Value = "נאמן"
$ws.Range("C42").Value = "עידו קרן"
$ws.Range("C43").Value = "תומר וקס"
$ws.Range("C44").Value = "תומר קדם"
$ws.Range("C45").Value = "אריאל ליבזון"
$ws.Range("C46").Value = "גל ארצי"
$ws.Range("C47").Value = "הילה"
$ws.Range("C48").Value = "ליאורה"
$ws.Range("C49").Value = "פסוול"
$ws.Range("C50").Value = "מנקר"
$ws.Range("C51").Value = "גרונר"
$ws.Range("C52").Value = "תמיר"
$ws.Range("C53").Value = "סטיב"
$ws.Range("C54").Value = "גלעד חננאל"
$ws.Range("C55").Value = "הוד"
$ws.Range("C56").Value = "זהר רטנר"
$ws.Range("C57").Value = "יואב סטרולוביץ'"
$ws.Range("C59").Value = "עמרי קונסטנטינו"
$ws.Range("C60").Value = "אביתר"
$ws.Range("C61").Value = "ווינטרויב"
$ws.Range("C62").Value = "גל נימצקי"
$ws.Range("C63").Value = "דנה"
$ws.Range("C64").Value = "יהל פורת"
$ws.Range("C65").Value = "שליו"
$ws.Range("C66").Value = "כהנא"
$ws.Range("C58").Value = "רומנו"

# Update the sheet view (pane freeze + active selection)
[void]$ws.Range("D56").Select()
